# Update "Overview" balance-sheet data: shift the four existing financial
# periods one column to the left (D<-E, E<-F, F<-G, G<-H) and populate the
# newly added right-most period (H) together with its refreshed header/
# publish-date labels - i.e. drop FY1396 and add FY1401, per the new
# "read_price" database refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial-period headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: report publish-date headers ---
$ws.Range("D9").Value = "1399-04-04 (8)"
$ws.Range("E9").Value = "1400-04-05 (11)"
$ws.Range("F9").Value = "1401-07-28 (13)"
$ws.Range("G9").Value = "1402-02-30 (10)"
$ws.Range("H9").Value = "1402-02-30"

# --- Rows 12-58: balance-sheet line items (skip blank section-header rows) ---
$rowData = @{
  12 = @(40952, 59352, 3913305, 3961691, 5020100)
  13 = @(9536875, 7579323, 46692929, 66218762, 70212245)
  14 = @(8763571, 8951384, 13653527, 13442429, 39593189)
  15 = @(1953206, 2735034, 4551340, 7820923, 7709363)
  16 = @(591108, 1172405, 1054175, 1108125, 1088950)
  17 = @(0, 0, 0, 0, 0)
  18 = @(20885712, 20497498, 69865276, 92551930, 123623847)
  19 = @(0, 0, 0, 0, 0)
  20 = @(3193323, 19110754, 20242843, 21544140, 43674784)
  21 = @(0, 0, 0, 0, 0)
  22 = @(880766, 742244, 627653, 1054682, 1958845)
  23 = @(26574, 26575, 26316, 25431, 32012)
  24 = @("-", "-", "-", "-", "-")
  25 = @(162257, 156884, 112367, 95529, 80318)
  26 = @(4262920, 20036457, 21009179, 22719782, 45745959)
  27 = @(25148632, 40533955, 90874455, 115271712, 169369806)
  29 = @(3915963, 6428293, 8993353, 25788558, 35818307)
  30 = @("-", "-", "-", "-", "-")
  31 = @(206207, 193753, 271054, 414709, 810377)
  32 = @(1072227, 1047242, 1004735, 5165117, 4598687)
  33 = @(40526, 212474, 71626, 508459, 109475)
  34 = @(1700000, 13500000, 15800000, 28700000, 57298482)
  35 = @(0, 0, 0, 0, 0)
  36 = @(0, 0, 0, 0, 0)
  37 = @(6934923, 21381762, 26140768, 60576843, 98635328)
  38 = @(0, 0, 0, 0, 0)
  39 = @("-", "-", "-", "-", "-")
  40 = @(0, 0, 0, 0, 0)
  41 = @(206123, 291724, 451952, 790659, 943756)
  42 = @(206123, 291724, 451952, 790659, 943756)
  43 = @(7141046, 21673486, 26592720, 61367502, 99579084)
  45 = @(950000, 950000, 950000, 950000, 30950000)
  46 = @(0, 0, 0, 0, 0)
  47 = @(0, 0, 0, 29777544, 0)
  48 = @(0, 0, -73711, -116032, -107152)
  49 = @(0, 0, 0, 33190, 23121)
  50 = @(100817, 100817, 100817, 100817, 1688779)
  51 = @(0, 0, 0, 0, 0)
  52 = @("-", "-", "-", "-", "-")
  53 = @(0, 0, 0, 0, 0)
  54 = @("-", "-", "-", "-", "-")
  55 = @(0, 0, 0, 0, 0)
  56 = @(16956769, 17809652, 63304629, 23158691, 37235974)
  57 = @(18007586, 18860469, 64281735, 53904210, 69790722)
  58 = @(25148632, 40533955, 90874455, 115271712, 169369806)
}

$cols = @("D","E","F","G","H")
foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$r").Value = $vals[$i]
  }
}

# --- Row-height tweaks carried over from the source edit (best-effort; the
# default-row-height / font-descent metadata itself is not settable through
# this object model). ---
$ws.Rows(2).RowHeight = 15.6
$ws.Rows(5).RowHeight = 40.8
$ws.Rows(6).RowHeight = 40.8
$ws.Rows(8).RowHeight = 15.6

